$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2699.4
$ws.Range("I19").Value = 800
$ws.Range("J19").Value = 3965.6667
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 3965.6667
$ws.Range("M19").Value = -625
$ws.Range("N19").Value = -4315.6667
$ws.Range("H32").Value = 2838.5293
$ws.Range("I32").Value = 2297.5
$ws.Range("J32").Value = 3319.4443
$ws.Range("K32").Value = 2297.5
$ws.Range("L32").Value = 3319.4443
$ws.Range("M32").Value = -1971.5
$ws.Range("N32").Value = -3971.4443
$ws.Range("H98").Value = 2284.8125
$ws.Range("I98").Value = 2145.2
$ws.Range("J98").Value = 4379
$ws.Range("K98").Value = 2145.2
$ws.Range("L98").Value = 4379
$ws.Range("M98").Value = -647.1999999999998
$ws.Range("N98").Value = -7375
$ws.Range("H106").Value = 26959.867
$ws.Range("I106").Value = 28581.545
$ws.Range("J106").Value = 22500.25
$ws.Range("K106").Value = 28581.545
$ws.Range("L106").Value = 22500.25
$ws.Range("M106").Value = -27950.545
$ws.Range("N106").Value = -23762.25
$ws.Range("H122").Value = 2284.8125
$ws.Range("I122").Value = 2145.2
$ws.Range("J122").Value = 4379
$ws.Range("K122").Value = 6435.599999999999
$ws.Range("L122").Value = 13137
$ws.Range("M122").Value = -3985.599999999999
$ws.Range("N122").Value = -18037
$ws.Range("H137").Value = 23816300
$ws.Range("I137").Value = 29419524
$ws.Range("J137").Value = 2604.75
$ws.Range("K137").Value = 88258572
$ws.Range("L137").Value = 7814.25
$ws.Range("M137").Value = -88256022
$ws.Range("N137").Value = -12914.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5220.8237
$ws.Range("I86").Value = 4937.3
$ws.Range("K86").Value = 4937.3
$ws.Range("M86").Value = -3814.3
$ws.Range("H89").Value = 5220.8237
$ws.Range("I89").Value = 4937.3
$ws.Range("K89").Value = 24686.5
$ws.Range("M89").Value = -19070.5
$ws.Range("H107").Value = 2149.524
$ws.Range("I107").Value = 2261.2942
$ws.Range("K107").Value = 2261.2942
$ws.Range("M107").Value = -341.2941999999998
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 253.4762
$ws.Range("I7").Value = 150.3
$ws.Range("J7").Value = 347.27274
$ws.Range("K7").Value = 150.3
$ws.Range("L7").Value = 347.27274
$ws.Range("M7").Value = -37.30000000000001
$ws.Range("N7").Value = -573.27274
$ws.Range("H14").Value = 5444.4443
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5340
$ws.Range("H16").Value = 2546.3684
$ws.Range("I16").Value = 1828.2858
$ws.Range("J16").Value = 2965.25
$ws.Range("K16").Value = 1828.2858
$ws.Range("L16").Value = 2965.25
$ws.Range("M16").Value = -1541.2858
$ws.Range("N16").Value = -3539.25
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H31").Value = 304985.62
$ws.Range("I31").Value = 7194.625
$ws.Range("J31").Value = 615724.0600000001
$ws.Range("K31").Value = 7194.625
$ws.Range("L31").Value = 615724.0600000001
$ws.Range("M31").Value = -6899.625
$ws.Range("N31").Value = -616314.0600000001
$ws.Range("H34").Value = 304985.62
$ws.Range("I34").Value = 7194.625
$ws.Range("J34").Value = 615724.0600000001
$ws.Range("K34").Value = 7194.625
$ws.Range("L34").Value = 615724.0600000001
$ws.Range("M34").Value = -6992.625
$ws.Range("N34").Value = -616128.0600000001
$ws.Range("H99").Value = 11350500
$ws.Range("I99").Value = 5266228.5
$ws.Range("K99").Value = 5266228.5
$ws.Range("M99").Value = -5264730.5
$ws.Range("H113").Value = 2546.3684
$ws.Range("I113").Value = 1828.2858
$ws.Range("J113").Value = 2965.25
$ws.Range("K113").Value = 1828.2858
$ws.Range("L113").Value = 2965.25
$ws.Range("M113").Value = 341.7141999999999
$ws.Range("N113").Value = -7305.25
$ws.Range("H122").Value = 2670.6667
$ws.Range("I122").Value = 2670.6667
$ws.Range("K122").Value = 8012.000100000001
$ws.Range("M122").Value = -5562.000100000001
$ws.Range("H126").Value = 11350500
$ws.Range("I126").Value = 5266228.5
$ws.Range("K126").Value = 15798685.5
$ws.Range("M126").Value = -15796215.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 750.8182
$ws.Range("I5").Value = 614.64703
$ws.Range("K5").Value = 1843.94109
$ws.Range("M5").Value = -1731.94109
$ws.Range("H70").Value = 125006250
$ws.Range("I70").Value = 125006250
$ws.Range("K70").Value = 375018750
$ws.Range("M70").Value = -375018435
$ws.Range("H73").Value = 125006250
$ws.Range("I73").Value = 125006250
$ws.Range("K73").Value = 375018750
$ws.Range("M73").Value = -375017658
$ws.Range("H131").Value = 1959.6493
$ws.Range("I131").Value = 1149.8889
$ws.Range("J131").Value = 2206.6948
$ws.Range("K131").Value = 3449.6667
$ws.Range("L131").Value = 6620.084400000001
$ws.Range("M131").Value = 1590.3333
$ws.Range("N131").Value = -16700.0844
$ws.Range("H135").Value = 750.8182
$ws.Range("I135").Value = 614.64703
$ws.Range("K135").Value = 5531.82327
$ws.Range("M135").Value = -2996.82327
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3073.3547
$ws.Range("I102").Value = 2728.2593
$ws.Range("K102").Value = 2728.2593
$ws.Range("M102").Value = -1106.2593
$ws.Range("H126").Value = 4028
$ws.Range("I126").Value = 4015.1
$ws.Range("K126").Value = 12045.3
$ws.Range("M126").Value = -9575.299999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9461.92
$ws.Range("I122").Value = 5015.7646
$ws.Range("K122").Value = 15047.2938
$ws.Range("M122").Value = -12597.2938
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 27945788
$ws.Range("I122").Value = 26790172
$ws.Range("K122").Value = 80370516
$ws.Range("M122").Value = -80368066
$ws.Range("H132").Value = 1425.8718
$ws.Range("I132").Value = 1443.6578
$ws.Range("K132").Value = 4330.9734
$ws.Range("M132").Value = -1800.9734
$ws.Range("H136").Value = 4680.1914
$ws.Range("J136").Value = 6316.6665
$ws.Range("L136").Value = 18949.9995
$ws.Range("N136").Value = -24049.9995
